$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextValue "D2" '58.049.37'
Set-TextValue "D3" '3.136.31'
Set-TextValue "E3" '  +1.44%  '
Set-TextValue "E4" '  -0.03%  '
Set-TextValue "D5" '535.12'
Set-TextValue "E5" '  +2.37%  '
Set-TextValue "E6" '  +2.11%  '
Set-TextValue "D7" '0.999'
Set-TextValue "E7" '  +0.01%  '
Set-TextValue "D8" '0.508'
Set-TextValue "E8" '  +11.41%  '
Set-TextValue "E9" '  -0.18%  '
Set-TextValue "B10" 'Cardano'
Set-TextValue "C10" 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
Set-TextValue "D10" '0.425'
Set-TextValue "E10" '  +6.76%  '
Set-TextValue "B11" 'Dogecoin'
Set-TextValue "C11" 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
Set-TextValue "D11" '0.109'
Set-TextValue "E11" '  +2.48%  '
Set-TextValue "E12" '  +3.44%  '
Set-TextValue "D13" '3.677.76'
Set-TextValue "E13" '  +1.41%  '
Set-TextValue "D14" '25.75'
Set-TextValue "E14" '  +1.59%  '
Set-TextValue "E15" '  +5.08%  '
Set-TextValue "D16" '58.112.77'
Set-TextValue "E16" '  +1.41%  '
Set-TextValue "D17" '6.26'
Set-TextValue "E17" '  +6.39%  '
Set-TextValue "D18" '3.137.53'
Set-TextValue "E18" '  +1.34%  '
Set-TextValue "E19" '  +4.19%  '
Set-TextValue "D20" '8.24'
Set-TextValue "E20" '  +4.92%  '
Set-TextValue "D21" '376.77'
Set-TextValue "E21" '  +7.97%  '
Set-TextValue "E22" '  -0.01%  '
Set-TextValue "D23" '5.74'
Set-TextValue "E23" '  -0.69%  '
Set-TextValue "D24" '70.20'
Set-TextValue "E24" '  +2.41%  '
Set-TextValue "E25" '  +3.49%  '
Set-TextValue "E26" '  +0.28%  '
Set-TextValue "D27" '1.00'
Set-TextValue "E27" '  +0.15%  '
Set-TextValue "B28" 'PEPE'
Set-TextValue "C28" 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
Set-TextValue "D28" '0.0₃0885'
Set-TextValue "E28" '  +1.84%  '
Set-TextValue "B29" 'InternetComputer(DFINITY)'
Set-TextValue "C29" 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextValue "D29" '7.95'
Set-TextValue "E29" '  +9.91%  '
Set-TextValue "E30" '  +5.12%  '
Set-TextValue "E31" '  +1.08%  '
Set-TextValue "D32" '21.76'
Set-TextValue "E32" '  +4.17%  '
Set-TextValue "D33" '5.18'
Set-TextValue "E33" '  +5.81%  '
Set-TextValue "D34" '1.17'
Set-TextValue "E34" '  +3.00%  '
Set-TextValue "D35" '161.53'
Set-TextValue "E35" '  +1.54%  '
Set-TextValue "E36" '  +4.31%  '
Set-TextValue "E37" '  +9.05%  '
Set-TextValue "D38" '25.54'
Set-TextValue "E38" '  -0.46%  '
Set-TextValue "E39" '  +5.39%  '
Set-TextValue "D40" '2.632.46'
Set-TextValue "E40" '  +9.60%  '
Set-TextValue "E41" '  +5.48%  '
Set-TextValue "E42" '  +2.43%  '
Set-TextValue "D43" '39.01'
Set-TextValue "E43" '  +6.28%  '
Set-TextValue "D44" '0.700'
Set-TextValue "E44" '  +0.60%  '
Set-TextValue "D45" '0.0273'
Set-TextValue "E45" '  +3.92%  '
Set-TextValue "E46" '  -0.13%  '
Set-TextValue "E47" '  +4.58%  '
Set-TextValue "D48" '0.978'
Set-TextValue "E48" '  +2.32%  '
Set-TextValue "D49" '0.0999'
Set-TextValue "E49" '  +10.18%  '
Set-TextValue "E50" '  +2.93%  '
Set-TextValue "E51" '  -1.61%  '
